$wb = $excel.ActiveWorkbook

# --- Sheet1: fix the boolean cell C7 to be a formula "=TRUE()" using the
#     default/general style instead of the dedicated TRUE/FALSE number format.
$ws1 = $wb.Worksheets.Item(1)
$c7 = $ws1.Range("C7")
$c7.Formula = "=TRUE()"
$c7.Style = "Normal"

# --- Sheet2: rename to "2ndsheet" and populate it with a small 3x3 table.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2ndsheet"

$ws2.Range("A1").Value = "column 1"
$ws2.Range("B1").Value = "column 2"
$ws2.Range("C1").Value = "last column"
$ws2.Range("A1:C1").Font.Bold = $true

$ws2.Range("A2").Value = "data 1,1"
$ws2.Range("B2").Value = "data 1,2"
$ws2.Range("C2").Value = "data 1,3"

$ws2.Range("A3").Value = "data 2,1"
$ws2.Range("B3").Value = "data 2,2"
$ws2.Range("C3").Value = "data 2,3"

# Sheet2 becomes the active/selected sheet, with B3 the selected cell.
$ws2.Activate()
$ws2.Range("B3").Select() | Out-Null
